# Convert the textual "N_1"/"N_2" values in columns A and B into real
# Excel numeric values (A => N.1, B => N), matching the "conversion for
# excel types" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: "1_1".."7_1" -> 1.1 .. 7.1 (row 5 stays blank, as before)
$ws.Range("A2").Value = 1.1
$ws.Range("A3").Value = 2.1
$ws.Range("A4").Value = 3.1
$ws.Range("A6").Value = 5.1
$ws.Range("A7").Value = 6.1
$ws.Range("A8").Value = 7.1

# Column B: "1_2".."7_2" -> 1 .. 7 (row 6 stays blank, as before)
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 2
$ws.Range("B4").Value = 3
$ws.Range("B5").Value = 4
$ws.Range("B7").Value = 6
$ws.Range("B8").Value = 7

# Update the active selection to the single cell B2.
$ws.Range("B2").Select()
